# Update TPM-derived NATMI metrics with new values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.820647333333334
$ws.Range("H2").Value = 26.461942
$ws.Range("I2").Value = 0.06415146660411865
$ws.Range("J2").Value = 0.06415146660411865
$ws.Range("M2").Value = 2.63379
$ws.Range("N2").Value = 7.90137
$ws.Range("O2").Value = 0.03319695559561149
$ws.Range("P2").Value = 0.03319695559561149
$ws.Range("Q2").Value = 23.23173274006
$ws.Range("R2").Value = 209.08559466054
$ws.Range("S2").Value = 0.00212963338825028
$ws.Range("T2").Value = 0.00212963338825028
$ws.Range("G3").Value = 8.820647333333334
$ws.Range("H3").Value = 26.461942
$ws.Range("I3").Value = 0.06415146660411865
$ws.Range("J3").Value = 0.06415146660411865
$ws.Range("O3").Value = 0.8262122860897556
$ws.Range("P3").Value = 0.8262122860897555
$ws.Range("Q3").Value = 578.195881899743
$ws.Range("R3").Value = 5203.762937097687
$ws.Range("S3").Value = 0.05300272987899948
$ws.Range("T3").Value = 0.05300272987899947
$ws.Range("G4").Value = 8.820647333333334
$ws.Range("H4").Value = 26.461942
$ws.Range("I4").Value = 0.06415146660411865
$ws.Range("J4").Value = 0.06415146660411865
$ws.Range("M4").Value = 10.26216366666667
$ws.Range("N4").Value = 30.786491
$ws.Range("O4").Value = 0.1293469075200494
$ws.Range("P4").Value = 0.1293469075200494
$ws.Range("Q4").Value = 90.51892658061357
$ws.Range("R4").Value = 814.6703392255221
$ws.Range("S4").Value = 0.008297793818118473
$ws.Range("T4").Value = 0.008297793818118473
$ws.Range("G5").Value = 8.820647333333334
$ws.Range("H5").Value = 26.461942
$ws.Range("I5").Value = 0.06415146660411865
$ws.Range("J5").Value = 0.06415146660411865
$ws.Range("M5").Value = 0.8920680000000001
$ws.Range("N5").Value = 2.676204
$ws.Range("O5").Value = 0.01124385079458346
$ws.Range("P5").Value = 0.01124385079458345
$ws.Range("Q5").Value = 7.868617225352001
$ws.Range("R5").Value = 70.817555028168
$ws.Range("S5").Value = 0.0007213095187504136
$ws.Range("T5").Value = 0.0007213095187504133
$ws.Range("I6").Value = 0.3979101621202897
$ws.Range("J6").Value = 0.3979101621202898
$ws.Range("M6").Value = 2.63379
$ws.Range("N6").Value = 7.90137
$ws.Range("O6").Value = 0.03319695559561149
$ws.Range("P6").Value = 0.03319695559561149
$ws.Range("Q6").Value = 144.09869376765
$ws.Range("R6").Value = 1296.88824390885
$ws.Range("S6").Value = 0.01320940598294983
$ws.Range("T6").Value = 0.01320940598294983
$ws.Range("I7").Value = 0.3979101621202897
$ws.Range("J7").Value = 0.3979101621202898
$ws.Range("O7").Value = 0.8262122860897556
$ws.Range("P7").Value = 0.8262122860897555
$ws.Range("S7").Value = 0.3287582647037499
$ws.Range("T7").Value = 0.3287582647037499
$ws.Range("I8").Value = 0.3979101621202897
$ws.Range("J8").Value = 0.3979101621202898
$ws.Range("M8").Value = 10.26216366666667
$ws.Range("N8").Value = 30.786491
$ws.Range("O8").Value = 0.1293469075200494
$ws.Range("P8").Value = 0.1293469075200494
$ws.Range("Q8").Value = 561.4587266245617
$ws.Range("R8").Value = 5053.128539621055
$ws.Range("S8").Value = 0.05146844894106099
$ws.Range("T8").Value = 0.051468448941061
$ws.Range("I9").Value = 0.3979101621202897
$ws.Range("J9").Value = 0.3979101621202898
$ws.Range("M9").Value = 0.8920680000000001
$ws.Range("N9").Value = 2.676204
$ws.Range("O9").Value = 0.01124385079458346
$ws.Range("P9").Value = 0.01124385079458345
$ws.Range("Q9").Value = 48.80640960438
$ws.Range("R9").Value = 439.25768643942
$ws.Range("S9").Value = 0.004474042492529052
$ws.Range("T9").Value = 0.004474042492529052
$ws.Range("G10").Value = 21.90816333333333
$ws.Range("H10").Value = 65.72449
$ws.Range("I10").Value = 0.1593353362087987
$ws.Range("J10").Value = 0.1593353362087987
$ws.Range("M10").Value = 2.63379
$ws.Range("N10").Value = 7.90137
$ws.Range("O10").Value = 0.03319695559561149
$ws.Range("P10").Value = 0.03319695559561149
$ws.Range("Q10").Value = 57.7015015057
$ws.Range("R10").Value = 519.3135135513
$ws.Range("S10").Value = 0.005289448080935317
$ws.Range("T10").Value = 0.005289448080935317
$ws.Range("G11").Value = 21.90816333333333
$ws.Range("H11").Value = 65.72449
$ws.Range("I11").Value = 0.1593353362087987
$ws.Range("J11").Value = 0.1593353362087987
$ws.Range("O11").Value = 0.8262122860897556
$ws.Range("P11").Value = 0.8262122860897555
$ws.Range("Q11").Value = 1436.086189666686
$ws.Range("R11").Value = 12924.77570700017
$ws.Range("S11").Value = 0.1316448123839514
$ws.Range("T11").Value = 0.1316448123839513
$ws.Range("G12").Value = 21.90816333333333
$ws.Range("H12").Value = 65.72449
$ws.Range("I12").Value = 0.1593353362087987
$ws.Range("J12").Value = 0.1593353362087987
$ws.Range("M12").Value = 10.26216366666667
$ws.Range("N12").Value = 30.786491
$ws.Range("O12").Value = 0.1293469075200494
$ws.Range("P12").Value = 0.1293469075200494
$ws.Range("Q12").Value = 224.8251577627323
$ws.Range("R12").Value = 2023.42641986459
$ws.Range("S12").Value = 0.02060953299727546
$ws.Range("T12").Value = 0.02060953299727546
$ws.Range("G13").Value = 21.90816333333333
$ws.Range("H13").Value = 65.72449
$ws.Range("I13").Value = 0.1593353362087987
$ws.Range("J13").Value = 0.1593353362087987
$ws.Range("M13").Value = 0.8920680000000001
$ws.Range("N13").Value = 2.676204
$ws.Range("O13").Value = 0.01124385079458346
$ws.Range("P13").Value = 0.01124385079458345
$ws.Range("Q13").Value = 19.54357144844
$ws.Range("R13").Value = 175.89214303596
$ws.Range("S13").Value = 0.001791542746636523
$ws.Range("T13").Value = 0.001791542746636523
$ws.Range("G14").Value = 52.056859
$ws.Range("H14").Value = 156.170577
$ws.Range("I14").Value = 0.3786030350667928
$ws.Range("J14").Value = 0.3786030350667929
$ws.Range("M14").Value = 2.63379
$ws.Range("N14").Value = 7.90137
$ws.Range("O14").Value = 0.03319695559561149
$ws.Range("P14").Value = 0.03319695559561149
$ws.Range("Q14").Value = 137.10683466561
$ws.Range("R14").Value = 1233.96151199049
$ws.Range("S14").Value = 0.01256846814347606
$ws.Range("T14").Value = 0.01256846814347606
$ws.Range("G15").Value = 52.056859
$ws.Range("H15").Value = 156.170577
$ws.Range("I15").Value = 0.3786030350667928
$ws.Range("J15").Value = 0.3786030350667929
$ws.Range("O15").Value = 0.8262122860897556
$ws.Range("P15").Value = 0.8262122860897555
$ws.Range("Q15").Value = 3412.341561904516
$ws.Range("R15").Value = 30711.07405714064
$ws.Range("S15").Value = 0.3128064791230548
$ws.Range("T15").Value = 0.3128064791230548
$ws.Range("G16").Value = 52.056859
$ws.Range("H16").Value = 156.170577
$ws.Range("I16").Value = 0.3786030350667928
$ws.Range("J16").Value = 0.3786030350667929
$ws.Range("M16").Value = 10.26216366666667
$ws.Range("N16").Value = 30.786491
$ws.Range("O16").Value = 0.1293469075200494
$ws.Range("P16").Value = 0.1293469075200494
$ws.Range("Q16").Value = 534.2160070305897
$ws.Range("R16").Value = 4807.944063275307
$ws.Range("S16").Value = 0.04897113176359448
$ws.Range("T16").Value = 0.04897113176359449
$ws.Range("G17").Value = 52.056859
$ws.Range("H17").Value = 156.170577
$ws.Range("I17").Value = 0.3786030350667928
$ws.Range("J17").Value = 0.3786030350667929
$ws.Range("M17").Value = 0.8920680000000001
$ws.Range("N17").Value = 2.676204
$ws.Range("O17").Value = 0.01124385079458346
$ws.Range("P17").Value = 0.01124385079458345
$ws.Range("Q17").Value = 46.438258094412
$ws.Range("R17").Value = 417.944322849708
$ws.Range("S17").Value = 0.004256956036667467
$ws.Range("T17").Value = 0.004256956036667467
